{"js": "// Polish translation update for Welcome.docx:\n// The SmartCash mining paragraph is trimmed so it now ends at\n// \"...quite some time.\" instead of\n// \"...quite some time, until Smartcash reaches a considerable market cap.\"\n\nconst oldTail =\n  \"quite some time, until Smartcash reaches a considerable market cap.\";\nconst newTail = \"quite some time.\";\n\nconst results = context.document.body.search(oldTail, {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newTail, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Polish translation update for Welcome.docx:\n# The SmartCash mining paragraph is trimmed so it now ends at\n# \"...quite some time.\" instead of\n# \"...quite some time, until Smartcash reaches a considerable market cap.\"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"quite some time, until Smartcash reaches a considerable market cap.\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"quite some time.\"\n$find.Forward = $true\n$find.Wrap = 0\n$find.MatchCase = $true\n$find.Execute([ref]$find.Text, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]0, [ref]$false, $find.Replacement.Text, 2)\n"}
